$d = $word.ActiveDocument

# 1. Update the Eclipse IDE product/version text.
$d.Content.Find.Execute("eclipse-SDK-4.5.1-win32-x86_64", $true, $false, $false, $false, $false,
                         $true, 1, $false, "eclipse-jee-mars-1-win32-x86_64.zip", 2) | Out-Null

# 2. Move the _GoBack bookmark from the end of the first paragraph to the
#    (now-empty) paragraph that immediately follows the Eclipse line.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*eclipse-jee-mars-1-win32-x86_64.zip*") {
        $target = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($target -ne $null) {
    $d.Bookmarks.Add("_GoBack", $target.Range) | Out-Null
}
